$wb = $excel.ActiveWorkbook

# --- Sheet index 1 ---
$ws = $wb.Worksheets.Item(1)
$refRow = 470

# Update existing overlapping rows in the block where the A value changed
$ws.Cells.Item(470, 1).Value = 45725.23111297454
$ws.Cells.Item(471, 1).Value = 45725.23111297454
$ws.Cells.Item(472, 1).Value = 45725.23111297454
$ws.Cells.Item(473, 1).Value = 45725.23111297454
$ws.Cells.Item(474, 1).Value = 45725.23111297454
$ws.Cells.Item(475, 1).Value = 45725.23113506944
$ws.Cells.Item(476, 1).Value = 45725.23113506944
$ws.Cells.Item(477, 1).Value = 45725.23113506944
$ws.Cells.Item(478, 1).Value = 45725.23113506944
$ws.Cells.Item(479, 1).Value = 45725.23113506944
$ws.Cells.Item(480, 1).Value = 45725.23115855324
$ws.Cells.Item(481, 1).Value = 45725.23115855324
$ws.Cells.Item(482, 1).Value = 45725.23115855324
$ws.Cells.Item(483, 1).Value = 45725.23115855324
$ws.Cells.Item(484, 1).Value = 45725.23115855324
$ws.Cells.Item(485, 1).Value = 45725.73125641204
$ws.Cells.Item(486, 1).Value = 45725.73125641204
$ws.Cells.Item(487, 1).Value = 45725.73125641204

# Append new rows at the end, copying B..I from the reference row (constant block content)
$refB = $ws.Cells.Item($refRow, 2).Value2
$refC = $ws.Cells.Item($refRow, 3).Value2
$refD = $ws.Cells.Item($refRow, 4).Value2
$refE = $ws.Cells.Item($refRow, 5).Value2
$refF = $ws.Cells.Item($refRow, 6).Value2
$refG = $ws.Cells.Item($refRow, 7).Value2
$refH = $ws.Cells.Item($refRow, 8).Value2
$refI = $ws.Cells.Item($refRow, 9).Value2
$refNumFmt = $ws.Cells.Item($refRow, 1).NumberFormat
$ws.Cells.Item(488, 1).Value = 45725.73127832176
$ws.Cells.Item(488, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(488, 2).Value = $refB
$ws.Cells.Item(488, 3).Value = $refC
$ws.Cells.Item(488, 4).Value = $refD
$ws.Cells.Item(488, 5).Value = $refE
$ws.Cells.Item(488, 6).Value = $refF
$ws.Cells.Item(488, 7).Value = $refG
$ws.Cells.Item(488, 8).Value = $refH
$ws.Cells.Item(488, 9).Value = $refI
$ws.Cells.Item(489, 1).Value = 45725.73127832176
$ws.Cells.Item(489, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(489, 2).Value = $refB
$ws.Cells.Item(489, 3).Value = $refC
$ws.Cells.Item(489, 4).Value = $refD
$ws.Cells.Item(489, 5).Value = $refE
$ws.Cells.Item(489, 6).Value = $refF
$ws.Cells.Item(489, 7).Value = $refG
$ws.Cells.Item(489, 8).Value = $refH
$ws.Cells.Item(489, 9).Value = $refI
$ws.Cells.Item(490, 1).Value = 45725.73127832176
$ws.Cells.Item(490, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(490, 2).Value = $refB
$ws.Cells.Item(490, 3).Value = $refC
$ws.Cells.Item(490, 4).Value = $refD
$ws.Cells.Item(490, 5).Value = $refE
$ws.Cells.Item(490, 6).Value = $refF
$ws.Cells.Item(490, 7).Value = $refG
$ws.Cells.Item(490, 8).Value = $refH
$ws.Cells.Item(490, 9).Value = $refI
$ws.Cells.Item(491, 1).Value = 45725.73130123843
$ws.Cells.Item(491, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(491, 2).Value = $refB
$ws.Cells.Item(491, 3).Value = $refC
$ws.Cells.Item(491, 4).Value = $refD
$ws.Cells.Item(491, 5).Value = $refE
$ws.Cells.Item(491, 6).Value = $refF
$ws.Cells.Item(491, 7).Value = $refG
$ws.Cells.Item(491, 8).Value = $refH
$ws.Cells.Item(491, 9).Value = $refI
$ws.Cells.Item(492, 1).Value = 45725.73130123843
$ws.Cells.Item(492, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(492, 2).Value = $refB
$ws.Cells.Item(492, 3).Value = $refC
$ws.Cells.Item(492, 4).Value = $refD
$ws.Cells.Item(492, 5).Value = $refE
$ws.Cells.Item(492, 6).Value = $refF
$ws.Cells.Item(492, 7).Value = $refG
$ws.Cells.Item(492, 8).Value = $refH
$ws.Cells.Item(492, 9).Value = $refI
$ws.Cells.Item(493, 1).Value = 45725.73130123843
$ws.Cells.Item(493, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(493, 2).Value = $refB
$ws.Cells.Item(493, 3).Value = $refC
$ws.Cells.Item(493, 4).Value = $refD
$ws.Cells.Item(493, 5).Value = $refE
$ws.Cells.Item(493, 6).Value = $refF
$ws.Cells.Item(493, 7).Value = $refG
$ws.Cells.Item(493, 8).Value = $refH
$ws.Cells.Item(493, 9).Value = $refI
$ws.Cells.Item(494, 1).Value = 45726.23139893518
$ws.Cells.Item(494, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(494, 2).Value = $refB
$ws.Cells.Item(494, 3).Value = $refC
$ws.Cells.Item(494, 4).Value = $refD
$ws.Cells.Item(494, 5).Value = $refE
$ws.Cells.Item(494, 6).Value = $refF
$ws.Cells.Item(494, 7).Value = $refG
$ws.Cells.Item(494, 8).Value = $refH
$ws.Cells.Item(494, 9).Value = $refI
$ws.Cells.Item(495, 1).Value = 45726.23142038195
$ws.Cells.Item(495, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(495, 2).Value = $refB
$ws.Cells.Item(495, 3).Value = $refC
$ws.Cells.Item(495, 4).Value = $refD
$ws.Cells.Item(495, 5).Value = $refE
$ws.Cells.Item(495, 6).Value = $refF
$ws.Cells.Item(495, 7).Value = $refG
$ws.Cells.Item(495, 8).Value = $refH
$ws.Cells.Item(495, 9).Value = $refI
$ws.Cells.Item(496, 1).Value = 45726.23144357639
$ws.Cells.Item(496, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(496, 2).Value = $refB
$ws.Cells.Item(496, 3).Value = $refC
$ws.Cells.Item(496, 4).Value = $refD
$ws.Cells.Item(496, 5).Value = $refE
$ws.Cells.Item(496, 6).Value = $refF
$ws.Cells.Item(496, 7).Value = $refG
$ws.Cells.Item(496, 8).Value = $refH
$ws.Cells.Item(496, 9).Value = $refI

# --- Sheet index 2 ---
$ws = $wb.Worksheets.Item(2)
$refRow = 484

# Update existing overlapping rows in the block where the A value changed
$ws.Cells.Item(484, 1).Value = 45725.07909302083
$ws.Cells.Item(485, 1).Value = 45725.07909302083
$ws.Cells.Item(486, 1).Value = 45725.07909302083
$ws.Cells.Item(487, 1).Value = 45725.07909302083
$ws.Cells.Item(488, 1).Value = 45725.07909302083
$ws.Cells.Item(489, 1).Value = 45725.07911518519
$ws.Cells.Item(490, 1).Value = 45725.07911518519
$ws.Cells.Item(491, 1).Value = 45725.07911518519
$ws.Cells.Item(492, 1).Value = 45725.07911518519
$ws.Cells.Item(493, 1).Value = 45725.07911518519
$ws.Cells.Item(494, 1).Value = 45725.07913833333
$ws.Cells.Item(495, 1).Value = 45725.07913833333
$ws.Cells.Item(496, 1).Value = 45725.07913833333
$ws.Cells.Item(497, 1).Value = 45725.07913833333
$ws.Cells.Item(498, 1).Value = 45725.07913833333
$ws.Cells.Item(499, 1).Value = 45725.57923533564
$ws.Cells.Item(500, 1).Value = 45725.57923533564
$ws.Cells.Item(501, 1).Value = 45725.57923533564

# Append new rows at the end, copying B..I from the reference row (constant block content)
$refB = $ws.Cells.Item($refRow, 2).Value2
$refC = $ws.Cells.Item($refRow, 3).Value2
$refD = $ws.Cells.Item($refRow, 4).Value2
$refE = $ws.Cells.Item($refRow, 5).Value2
$refF = $ws.Cells.Item($refRow, 6).Value2
$refG = $ws.Cells.Item($refRow, 7).Value2
$refH = $ws.Cells.Item($refRow, 8).Value2
$refI = $ws.Cells.Item($refRow, 9).Value2
$refNumFmt = $ws.Cells.Item($refRow, 1).NumberFormat
$ws.Cells.Item(502, 1).Value = 45725.57925716435
$ws.Cells.Item(502, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(502, 2).Value = $refB
$ws.Cells.Item(502, 3).Value = $refC
$ws.Cells.Item(502, 4).Value = $refD
$ws.Cells.Item(502, 5).Value = $refE
$ws.Cells.Item(502, 6).Value = $refF
$ws.Cells.Item(502, 7).Value = $refG
$ws.Cells.Item(502, 8).Value = $refH
$ws.Cells.Item(502, 9).Value = $refI
$ws.Cells.Item(503, 1).Value = 45725.57925716435
$ws.Cells.Item(503, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(503, 2).Value = $refB
$ws.Cells.Item(503, 3).Value = $refC
$ws.Cells.Item(503, 4).Value = $refD
$ws.Cells.Item(503, 5).Value = $refE
$ws.Cells.Item(503, 6).Value = $refF
$ws.Cells.Item(503, 7).Value = $refG
$ws.Cells.Item(503, 8).Value = $refH
$ws.Cells.Item(503, 9).Value = $refI
$ws.Cells.Item(504, 1).Value = 45725.57925716435
$ws.Cells.Item(504, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(504, 2).Value = $refB
$ws.Cells.Item(504, 3).Value = $refC
$ws.Cells.Item(504, 4).Value = $refD
$ws.Cells.Item(504, 5).Value = $refE
$ws.Cells.Item(504, 6).Value = $refF
$ws.Cells.Item(504, 7).Value = $refG
$ws.Cells.Item(504, 8).Value = $refH
$ws.Cells.Item(504, 9).Value = $refI
$ws.Cells.Item(505, 1).Value = 45725.57928042824
$ws.Cells.Item(505, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(505, 2).Value = $refB
$ws.Cells.Item(505, 3).Value = $refC
$ws.Cells.Item(505, 4).Value = $refD
$ws.Cells.Item(505, 5).Value = $refE
$ws.Cells.Item(505, 6).Value = $refF
$ws.Cells.Item(505, 7).Value = $refG
$ws.Cells.Item(505, 8).Value = $refH
$ws.Cells.Item(505, 9).Value = $refI
$ws.Cells.Item(506, 1).Value = 45725.57928042824
$ws.Cells.Item(506, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(506, 2).Value = $refB
$ws.Cells.Item(506, 3).Value = $refC
$ws.Cells.Item(506, 4).Value = $refD
$ws.Cells.Item(506, 5).Value = $refE
$ws.Cells.Item(506, 6).Value = $refF
$ws.Cells.Item(506, 7).Value = $refG
$ws.Cells.Item(506, 8).Value = $refH
$ws.Cells.Item(506, 9).Value = $refI
$ws.Cells.Item(507, 1).Value = 45725.57928042824
$ws.Cells.Item(507, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(507, 2).Value = $refB
$ws.Cells.Item(507, 3).Value = $refC
$ws.Cells.Item(507, 4).Value = $refD
$ws.Cells.Item(507, 5).Value = $refE
$ws.Cells.Item(507, 6).Value = $refF
$ws.Cells.Item(507, 7).Value = $refG
$ws.Cells.Item(507, 8).Value = $refH
$ws.Cells.Item(507, 9).Value = $refI
$ws.Cells.Item(508, 1).Value = 45726.07937777778
$ws.Cells.Item(508, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(508, 2).Value = $refB
$ws.Cells.Item(508, 3).Value = $refC
$ws.Cells.Item(508, 4).Value = $refD
$ws.Cells.Item(508, 5).Value = $refE
$ws.Cells.Item(508, 6).Value = $refF
$ws.Cells.Item(508, 7).Value = $refG
$ws.Cells.Item(508, 8).Value = $refH
$ws.Cells.Item(508, 9).Value = $refI
$ws.Cells.Item(509, 1).Value = 45726.07939922454
$ws.Cells.Item(509, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(509, 2).Value = $refB
$ws.Cells.Item(509, 3).Value = $refC
$ws.Cells.Item(509, 4).Value = $refD
$ws.Cells.Item(509, 5).Value = $refE
$ws.Cells.Item(509, 6).Value = $refF
$ws.Cells.Item(509, 7).Value = $refG
$ws.Cells.Item(509, 8).Value = $refH
$ws.Cells.Item(509, 9).Value = $refI
$ws.Cells.Item(510, 1).Value = 45726.07942256945
$ws.Cells.Item(510, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(510, 2).Value = $refB
$ws.Cells.Item(510, 3).Value = $refC
$ws.Cells.Item(510, 4).Value = $refD
$ws.Cells.Item(510, 5).Value = $refE
$ws.Cells.Item(510, 6).Value = $refF
$ws.Cells.Item(510, 7).Value = $refG
$ws.Cells.Item(510, 8).Value = $refH
$ws.Cells.Item(510, 9).Value = $refI

# --- Sheet index 4 ---
$ws = $wb.Worksheets.Item(4)
$refRow = 488

# Update existing overlapping rows in the block where the A value changed
$ws.Cells.Item(488, 1).Value = 45725.22919952546
$ws.Cells.Item(489, 1).Value = 45725.22919952546
$ws.Cells.Item(490, 1).Value = 45725.22919952546
$ws.Cells.Item(491, 1).Value = 45725.22919952546
$ws.Cells.Item(492, 1).Value = 45725.22919952546
$ws.Cells.Item(493, 1).Value = 45725.22922125
$ws.Cells.Item(494, 1).Value = 45725.22922125
$ws.Cells.Item(495, 1).Value = 45725.22922125
$ws.Cells.Item(496, 1).Value = 45725.22922125
$ws.Cells.Item(497, 1).Value = 45725.22922125
$ws.Cells.Item(498, 1).Value = 45725.22924497685
$ws.Cells.Item(499, 1).Value = 45725.22924497685
$ws.Cells.Item(500, 1).Value = 45725.22924497685
$ws.Cells.Item(501, 1).Value = 45725.22924497685
$ws.Cells.Item(502, 1).Value = 45725.22924497685
$ws.Cells.Item(503, 1).Value = 45725.7293421412
$ws.Cells.Item(504, 1).Value = 45725.7293421412
$ws.Cells.Item(505, 1).Value = 45725.7293421412

# Append new rows at the end, copying B..I from the reference row (constant block content)
$refB = $ws.Cells.Item($refRow, 2).Value2
$refC = $ws.Cells.Item($refRow, 3).Value2
$refD = $ws.Cells.Item($refRow, 4).Value2
$refE = $ws.Cells.Item($refRow, 5).Value2
$refF = $ws.Cells.Item($refRow, 6).Value2
$refG = $ws.Cells.Item($refRow, 7).Value2
$refH = $ws.Cells.Item($refRow, 8).Value2
$refI = $ws.Cells.Item($refRow, 9).Value2
$refNumFmt = $ws.Cells.Item($refRow, 1).NumberFormat
$ws.Cells.Item(506, 1).Value = 45725.72936453704
$ws.Cells.Item(506, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(506, 2).Value = $refB
$ws.Cells.Item(506, 3).Value = $refC
$ws.Cells.Item(506, 4).Value = $refD
$ws.Cells.Item(506, 5).Value = $refE
$ws.Cells.Item(506, 6).Value = $refF
$ws.Cells.Item(506, 7).Value = $refG
$ws.Cells.Item(506, 8).Value = $refH
$ws.Cells.Item(506, 9).Value = $refI
$ws.Cells.Item(507, 1).Value = 45725.72936453704
$ws.Cells.Item(507, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(507, 2).Value = $refB
$ws.Cells.Item(507, 3).Value = $refC
$ws.Cells.Item(507, 4).Value = $refD
$ws.Cells.Item(507, 5).Value = $refE
$ws.Cells.Item(507, 6).Value = $refF
$ws.Cells.Item(507, 7).Value = $refG
$ws.Cells.Item(507, 8).Value = $refH
$ws.Cells.Item(507, 9).Value = $refI
$ws.Cells.Item(508, 1).Value = 45725.72936453704
$ws.Cells.Item(508, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(508, 2).Value = $refB
$ws.Cells.Item(508, 3).Value = $refC
$ws.Cells.Item(508, 4).Value = $refD
$ws.Cells.Item(508, 5).Value = $refE
$ws.Cells.Item(508, 6).Value = $refF
$ws.Cells.Item(508, 7).Value = $refG
$ws.Cells.Item(508, 8).Value = $refH
$ws.Cells.Item(508, 9).Value = $refI
$ws.Cells.Item(509, 1).Value = 45725.7293875
$ws.Cells.Item(509, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(509, 2).Value = $refB
$ws.Cells.Item(509, 3).Value = $refC
$ws.Cells.Item(509, 4).Value = $refD
$ws.Cells.Item(509, 5).Value = $refE
$ws.Cells.Item(509, 6).Value = $refF
$ws.Cells.Item(509, 7).Value = $refG
$ws.Cells.Item(509, 8).Value = $refH
$ws.Cells.Item(509, 9).Value = $refI
$ws.Cells.Item(510, 1).Value = 45725.7293875
$ws.Cells.Item(510, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(510, 2).Value = $refB
$ws.Cells.Item(510, 3).Value = $refC
$ws.Cells.Item(510, 4).Value = $refD
$ws.Cells.Item(510, 5).Value = $refE
$ws.Cells.Item(510, 6).Value = $refF
$ws.Cells.Item(510, 7).Value = $refG
$ws.Cells.Item(510, 8).Value = $refH
$ws.Cells.Item(510, 9).Value = $refI
$ws.Cells.Item(511, 1).Value = 45725.7293875
$ws.Cells.Item(511, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(511, 2).Value = $refB
$ws.Cells.Item(511, 3).Value = $refC
$ws.Cells.Item(511, 4).Value = $refD
$ws.Cells.Item(511, 5).Value = $refE
$ws.Cells.Item(511, 6).Value = $refF
$ws.Cells.Item(511, 7).Value = $refG
$ws.Cells.Item(511, 8).Value = $refH
$ws.Cells.Item(511, 9).Value = $refI
$ws.Cells.Item(512, 1).Value = 45726.2294844213
$ws.Cells.Item(512, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(512, 2).Value = $refB
$ws.Cells.Item(512, 3).Value = $refC
$ws.Cells.Item(512, 4).Value = $refD
$ws.Cells.Item(512, 5).Value = $refE
$ws.Cells.Item(512, 6).Value = $refF
$ws.Cells.Item(512, 7).Value = $refG
$ws.Cells.Item(512, 8).Value = $refH
$ws.Cells.Item(512, 9).Value = $refI
$ws.Cells.Item(513, 1).Value = 45726.22950657408
$ws.Cells.Item(513, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(513, 2).Value = $refB
$ws.Cells.Item(513, 3).Value = $refC
$ws.Cells.Item(513, 4).Value = $refD
$ws.Cells.Item(513, 5).Value = $refE
$ws.Cells.Item(513, 6).Value = $refF
$ws.Cells.Item(513, 7).Value = $refG
$ws.Cells.Item(513, 8).Value = $refH
$ws.Cells.Item(513, 9).Value = $refI
$ws.Cells.Item(514, 1).Value = 45726.22953034722
$ws.Cells.Item(514, 1).NumberFormat = $refNumFmt
$ws.Cells.Item(514, 2).Value = $refB
$ws.Cells.Item(514, 3).Value = $refC
$ws.Cells.Item(514, 4).Value = $refD
$ws.Cells.Item(514, 5).Value = $refE
$ws.Cells.Item(514, 6).Value = $refF
$ws.Cells.Item(514, 7).Value = $refG
$ws.Cells.Item(514, 8).Value = $refH
$ws.Cells.Item(514, 9).Value = $refI
